# Realestate Update resale numbers 2023-07-02 09:53
# Appends a new data row (row 99) to the CityResaleNum sheet, mirroring the
# layout of the existing rows (Date, Time, Weekday, Week as text; the
# remaining metric columns as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 99

# Values that must land in the sheet as literal TEXT (matching columns A-D
# of every other row), even though some of them look numeric/date-like.
$textValues = @{
    "A" = "2023-07-02"
    "B" = "09:48:22"
    "C" = "Sunday"
    "D" = "27"
}

# Values that are genuine numbers (columns E-T).
$numberValues = @{
    "E" = 123559
    "F" = 135013
    "G" = 161299
    "H" = 131594
    "I" = 175703
    "J" = 113365
    "K" = 204927
    "L" = 222598
    "M" = 174732
    "N" = 103471
    "O" = 38883
    "P" = 32901
    "Q" = 51985
    "R" = -1
    "S" = 35836
    "T" = -1
}

# Scratch cell well outside the sheet's real data (A1:T98) used to stage a
# text-formula result so it can be pasted as a *value* into the target cell.
# Routing through a formula ("=\"2023-07-02\"") guarantees the pasted result
# keeps a text type instead of Excel's usual "looks like a number/date" auto
# conversion that a direct .Value assignment would trigger.
$helper = $ws.Range("AA1")

foreach ($col in "A", "B", "C", "D") {
    $text = $textValues[$col]
    $target = $ws.Range($col + $newRow)

    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $target.PasteSpecial(-4163)
}

$excel.CutCopyMode = $false
$helper.Clear()

foreach ($col in "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T") {
    $ws.Range($col + $newRow).Value = $numberValues[$col]
}
